# Applies scheduled-runner price/profit updates to the Cactuar_Profits workbook.
# Each worksheet (named by job abbreviation) gets refreshed currentAveragePrice /
# LevePrice / LeveProfit figures in columns H-N for specific leve rows.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2066.6667
$ws.Range("J29").Value = 1750
$ws.Range("L29").Value = 5250
$ws.Range("N29").Value = -5812
$ws.Range("H33").Value = 648.4666999999999
$ws.Range("I33").Value = 680.5
$ws.Range("K33").Value = 680.5
$ws.Range("M33").Value = -451.5
$ws.Range("H40").Value = 13909183
$ws.Range("I40").Value = 11985.048
$ws.Range("J40").Value = 33365260
$ws.Range("K40").Value = 11985.048
$ws.Range("L40").Value = 33365260
$ws.Range("M40").Value = -11810.048
$ws.Range("N40").Value = -33365610
$ws.Range("H43").Value = 3080510.5
$ws.Range("I43").Value = 5131104.5
$ws.Range("K43").Value = 5131104.5
$ws.Range("M43").Value = -5131035.5
$ws.Range("H98").Value = 1447.7646
$ws.Range("I98").Value = 1514.5625
$ws.Range("K98").Value = 1514.5625
$ws.Range("M98").Value = -16.5625
$ws.Range("H100").Value = 1556
$ws.Range("I100").Value = 1556
$ws.Range("K100").Value = 1556
$ws.Range("M100").Value = -1015
$ws.Range("H101").Value = 776.8
$ws.Range("J101").Value = 858.5
$ws.Range("L101").Value = 2575.5
$ws.Range("N101").Value = -5819.5
$ws.Range("H107").Value = 213.8
$ws.Range("I107").Value = 163.33333
$ws.Range("J107").Value = 289.5
$ws.Range("K107").Value = 163.33333
$ws.Range("L107").Value = 289.5
$ws.Range("M107").Value = 1756.66667
$ws.Range("N107").Value = -4129.5
$ws.Range("H112").Value = 3001.432
$ws.Range("J112").Value = 3121.2927
$ws.Range("L112").Value = 9363.8781
$ws.Range("N112").Value = -11579.8781
$ws.Range("H113").Value = 2971.76
$ws.Range("I113").Value = 1828.5
$ws.Range("K113").Value = 1828.5
$ws.Range("M113").Value = 1425.5
$ws.Range("H118").Value = 1636
$ws.Range("I118").Value = 1636
$ws.Range("K118").Value = 4908
$ws.Range("M118").Value = -3251
$ws.Range("H122").Value = 1447.7646
$ws.Range("I122").Value = 1514.5625
$ws.Range("K122").Value = 4543.6875
$ws.Range("M122").Value = -2093.6875
$ws.Range("H132").Value = 153423.5
$ws.Range("I132").Value = 413298
$ws.Range("K132").Value = 1239894
$ws.Range("M132").Value = -1237364
$ws.Range("H137").Value = 8378124.5
$ws.Range("I137").Value = 527854.6
$ws.Range("K137").Value = 1583563.8
$ws.Range("M137").Value = -1581013.8
$ws.Range("H138").Value = 5475.7256
$ws.Range("J138").Value = 5922.8887
$ws.Range("L138").Value = 17768.6661
$ws.Range("N138").Value = -28048.6661

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4084.5652
$ws.Range("I32").Value = 2142.2373
$ws.Range("J32").Value = 15544.3
$ws.Range("K32").Value = 2142.2373
$ws.Range("L32").Value = 15544.3
$ws.Range("M32").Value = -1855.2373
$ws.Range("N32").Value = -16118.3
$ws.Range("H74").Value = 27780366
$ws.Range("I74").Value = 62500950
$ws.Range("K74").Value = 62500950
$ws.Range("M74").Value = -62500076
$ws.Range("H77").Value = 27780366
$ws.Range("I77").Value = 62500950
$ws.Range("K77").Value = 312504750
$ws.Range("M77").Value = -312500382
$ws.Range("H97").Value = 437.70587
$ws.Range("I97").Value = 430.46155
$ws.Range("K97").Value = 430.46155
$ws.Range("M97").Value = 65.53845000000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H64").Value = 4630137.5
$ws.Range("I64").Value = 6944793.5
$ws.Range("J64").Value = 825.3333
$ws.Range("K64").Value = 6944793.5
$ws.Range("L64").Value = 825.3333
$ws.Range("M64").Value = -6944568.5
$ws.Range("N64").Value = -1275.3333
$ws.Range("H67").Value = 4630137.5
$ws.Range("I67").Value = 6944793.5
$ws.Range("J67").Value = 825.3333
$ws.Range("K67").Value = 6944793.5
$ws.Range("L67").Value = 825.3333
$ws.Range("M67").Value = -6944013.5
$ws.Range("N67").Value = -2385.3333
$ws.Range("H134").Value = 2670.48
$ws.Range("I134").Value = 932.6111
$ws.Range("K134").Value = 2797.8333
$ws.Range("M134").Value = -262.8332999999998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1699.0625
$ws.Range("I16").Value = 1789.6364
$ws.Range("K16").Value = 1789.6364
$ws.Range("M16").Value = -1502.6364
$ws.Range("H31").Value = 5196.3193
$ws.Range("I31").Value = 2370.9092
$ws.Range("K31").Value = 2370.9092
$ws.Range("M31").Value = -2075.9092
$ws.Range("H34").Value = 5196.3193
$ws.Range("I34").Value = 2370.9092
$ws.Range("K34").Value = 2370.9092
$ws.Range("M34").Value = -2168.9092
$ws.Range("H113").Value = 1699.0625
$ws.Range("I113").Value = 1789.6364
$ws.Range("K113").Value = 1789.6364
$ws.Range("M113").Value = 380.3635999999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 274.66666
$ws.Range("I14").Value = 274.66666
$ws.Range("K14").Value = 823.9999799999999
$ws.Range("M14").Value = -650.9999799999999
$ws.Range("H17").Value = 4999
$ws.Range("J17").Value = 4999
$ws.Range("L17").Value = 14997
$ws.Range("N17").Value = -15335
$ws.Range("H46").Value = 3665.7778
$ws.Range("J46").Value = 4998.6665
$ws.Range("L46").Value = 14995.9995
$ws.Range("N46").Value = -15177.9995
$ws.Range("H56").Value = 6715.5557
$ws.Range("I56").Value = 6715.5557
$ws.Range("K56").Value = 6715.5557
$ws.Range("M56").Value = -6185.5557
$ws.Range("H112").Value = 10249.833
$ws.Range("I112").Value = 2874.75
$ws.Range("K112").Value = 8624.25
$ws.Range("M112").Value = -7516.25
$ws.Range("H113").Value = 676.1667
$ws.Range("I113").Value = 589.625
$ws.Range("K113").Value = 1768.875
$ws.Range("M113").Value = 401.125

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 471.75
$ws.Range("I97").Value = 463.19232
$ws.Range("K97").Value = 463.19232
$ws.Range("M97").Value = 32.80768
$ws.Range("H113").Value = 1304
$ws.Range("J113").Value = 1317.2
$ws.Range("L113").Value = 1317.2
$ws.Range("N113").Value = -5657.2
$ws.Range("H132").Value = 115412.72
$ws.Range("I132").Value = 157503.47
$ws.Range("J132").Value = 5976.8
$ws.Range("K132").Value = 472510.41
$ws.Range("L132").Value = 17930.4
$ws.Range("M132").Value = -469980.41
$ws.Range("N132").Value = -22990.4

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19609854
$ws.Range("I40").Value = 2093.3076
$ws.Range("K40").Value = 2093.3076
$ws.Range("M40").Value = -1957.3076
$ws.Range("H93").Value = 1649.4
$ws.Range("I93").Value = 1374.25
$ws.Range("K93").Value = 1374.25
$ws.Range("M93").Value = -126.25
$ws.Range("H122").Value = 57147224
$ws.Range("I122").Value = 71432610
$ws.Range("K122").Value = 214297830
$ws.Range("M122").Value = -214295380
$ws.Range("H132").Value = 4233.6665
$ws.Range("J132").Value = 4989
$ws.Range("L132").Value = 14967
$ws.Range("N132").Value = -20027
$ws.Range("H136").Value = 4516.2856
$ws.Range("I136").Value = 2423.6155
$ws.Range("J136").Value = 7916.875
$ws.Range("K136").Value = 7270.8465
$ws.Range("L136").Value = 23750.625
$ws.Range("M136").Value = -4720.8465
$ws.Range("N136").Value = -28850.625

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 94845.5
$ws.Range("I57").Value = 69696
$ws.Range("K57").Value = 69696
$ws.Range("M57").Value = -68942
$ws.Range("H64").Value = 94995
$ws.Range("J64").Value = 94995
$ws.Range("L64").Value = 94995
$ws.Range("N64").Value = -95491
$ws.Range("H67").Value = 94995
$ws.Range("J67").Value = 94995
$ws.Range("L67").Value = 94995
$ws.Range("N67").Value = -96711
$ws.Range("H100").Value = 823151
$ws.Range("I100").Value = 1327869.4
$ws.Range("J100").Value = 2983.625
$ws.Range("K100").Value = 2655738.8
$ws.Range("L100").Value = 5967.25
$ws.Range("M100").Value = -2655197.8
$ws.Range("N100").Value = -7049.25
$ws.Range("H132").Value = 111113510
$ws.Range("I132").Value = 27778778
$ws.Range("K132").Value = 83336334
$ws.Range("M132").Value = -83333804
